$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.626.00"
$ws.Range("E2").Value = "  +6.52%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.042.45"
$ws.Range("E3").Value = "  +3.34%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.90"
$ws.Range("E5").Value = "  +5.13%  "

$ws.Range("E6").Value = "  +2.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "66.42"
$ws.Range("E7").Value = "  +18.89%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  +6.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.44"
$ws.Range("E10").Value = "  +0.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0753"
$ws.Range("E11").Value = "  +4.35%  "

$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.907"
$ws.Range("E13").Value = "  +2.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.10"
$ws.Range("E14").Value = "  +6.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.343.20"
$ws.Range("E15").Value = "  +3.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.62"
$ws.Range("E16").Value = "  +7.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.83"
$ws.Range("E17").Value = "  +23.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.058.06"
$ws.Range("E18").Value = "  +4.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.433.14"
$ws.Range("E19").Value = "  +6.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.35"
$ws.Range("E20").Value = "  +5.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0874"
$ws.Range("E21").Value = "  +5.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.35"
$ws.Range("E22").Value = "  +7.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.31"
$ws.Range("E23").Value = "  +2.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.69"
$ws.Range("E24").Value = "  +20.64%  "

$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("E26").Value = "  +5.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.56"
$ws.Range("E27").Value = "  +6.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.35"
$ws.Range("E28").Value = "  +1.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.87"
$ws.Range("E29").Value = "  +2.88%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.23"
$ws.Range("E30").Value = "  +10.63%  "

$ws.Range("E31").Value = "  +3.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.21"
$ws.Range("E32").Value = "  +7.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.111"
$ws.Range("E33").Value = "  +24.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.73"
$ws.Range("E34").Value = "  +12.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0612"
$ws.Range("E35").Value = "  +5.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("E36").Value = "  +9.11%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.10"
$ws.Range("E38").Value = "  +26.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.82"
$ws.Range("E39").Value = "  +1.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.103"
$ws.Range("E40").Value = "  +17.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.23"
$ws.Range("E41").Value = "  +4.49%  "

$ws.Range("E42").Value = "  +4.74%  "

$ws.Range("E43").Value = "  +6.19%  "

$ws.Range("E44").Value = "  +22.28%  "

$ws.Range("E45").Value = "  +6.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.07"
$ws.Range("E46").Value = "  +9.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.95"
$ws.Range("E47").Value = "  +10.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.21"
$ws.Range("E48").Value = "  +5.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.426.04"
$ws.Range("E49").Value = "  +6.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.93"
$ws.Range("E50").Value = "  +2.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.25"
$ws.Range("E51").Value = "  +4.82%  "
